# Christmas Tree doc: add a "Meta description" paragraph right after the
# title heading, drop the duplicate bold title paragraph near the end, and
# replace the italic blurb there with the image-generation prompt text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the new "Meta description" paragraph right after paragraph 1
#    (the Heading1 title). It needs three runs: an empty leading run, a
#    bold "Meta description" run, and a normal run with the rest of the
#    text -- matching the structure used elsewhere in this document.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">: Experience festive fun with Christmas Tree. Play Christmas Tree for free with excellent graphics, exciting features, and top payouts.</w:t></w:r>' + `
  '</w:p>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Near the end of the document, the paragraph that duplicated the bold
#    title ("Play Christmas Tree Free| Exciting Features & Festive
#    Theme") is removed entirely -- its content now lives in the meta
#    description paragraph above.
# ---------------------------------------------------------------------
$oldTitlePara = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Play Christmas Tree Free| Exciting Features & Festive Theme") {
        $oldTitlePara = $p
        break
    }
}
if ($oldTitlePara -ne $null) {
    $oldTitlePara.Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------
# 3) The final paragraph (the italic blurb) keeps its formatting, but its
#    text is replaced with the AI image-generation prompt.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute(
    "Experience festive fun with Christmas Tree. Play Christmas Tree for free with excellent graphics, exciting features, and top payouts.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Prompt: Please create a feature image for the Christmas Tree online slot game. The image should be in cartoon style and feature a happy Maya warrior wearing glasses. The image should also incorporate Christmas elements such as snowflakes, gift boxes, and a decorated Christmas tree. The overall theme should be fun and festive, showcasing the excitement of the holiday season and the thrill of playing this slot game.",
    2
) | Out-Null
